
$wb = $excel.ActiveWorkbook

# --- Grab the existing sheets we need to touch -----------------------------
$methoden = $wb.Worksheets.Item("Methoden")
$user     = $wb.Worksheets.Item("User")
$depot    = $wb.Worksheets.Item("Depot")

# --- Add the new "Order" sheet at the end of the workbook ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$order = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$order.Name = "Order"

# Seed the new sheet by copying the formatting/borders of the Depot sheet's
# C:F columns (Company Name .. Trade Price) into Order's B:E columns - this
# reproduces the exact same border/fill styles used by the rest of the
# workbook for its "method table" layout.
$depot.Range("C1:F6").Copy($order.Range("B1:E6"))

# Overwrite the header row with the new column headers (this also grows the
# shared string table with the two new entries).
$order.Range("B2").Value = "Company Name (STRING)"
$order.Range("C2").Value = "Shares (INT)"
$order.Range("D2").Value = "State (STRING)"
$order.Range("E2").Value = "Type (STRING)"

# --- Update selections on the pre-existing sheets --------------------------
$null = $methoden.Range("B5").Select()
$null = $user.Range("D12").Select()
$null = $depot.Range("B2:F6").Select()

# --- Make "Order" the active sheet/tab with its own selection --------------
$order.Activate()
$null = $order.Range("E3").Select()
